$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workblocks")

# Insert 2 new rows after row 4 (before row 5), shifting the existing
# workblock rows (old 5-16) down to 7-18 to make room for a new
# "RecoverApps" workblock entry.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# New "wbCloseAppsRecover" workblock rows - set the label (column A) cells
# first so the new unique strings land in the shared string table in the
# same order the original author's edit produced them.
$ws.Range("A5").Value = "wbCloseAppsRecover_Type"
$ws.Range("A6").Value = "wbCloseAppsRecover_SuppressSuccessful"

# Shorten the "Value" column (B) text for every workblock: the
# "Task1, <Layer>, " prefix is dropped, keeping only the short name.
$ws.Range("B3").Value = "Init"
$ws.Range("B5").Value = "RecoverApps"
$ws.Range("B7").Value = "GetData"
$ws.Range("B9").Value = "Process"
$ws.Range("B11").Value = "Next"
$ws.Range("B13").Value = "CloseApps"
$ws.Range("B15").Value = "InitApps"
$ws.Range("B17").Value = "ProcessApps"

# Fill in the remaining cells for the new RecoverApps rows.
$ws.Range("C5").Value = "Name of Workblock"
$ws.Range("B6").Value = $true
$ws.Range("C6").Value = "Do not log successful executions of wb"

# The trailing two rows (wbProcess_Type / wbProcess_SuppressSuccessful)
# lost their explicit cell formatting in the source edit.
$ws.Range("B17").ClearFormats()
$ws.Range("C17").ClearFormats()
$ws.Range("B18").ClearFormats()
$ws.Range("C18").ClearFormats()

# Make "Workblocks" the active sheet/selection (was "Constants" before).
$ws.Activate()
$ws.Range("A3:C18").Select()
